# Auto-update data + news
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (A1:I1) loses its bold / bordered / centered formatting and
# reverts to the default "Normal" style.
$ws.Range("A1:I1").ClearFormats()

# Refresh the latest ICSA (Initial Jobless Claims) data row with the
# newly-observed reading.
$ws.Range("E9").Value = 200000
$ws.Range("G9").Value = 364147.5095785441
$ws.Range("H9").Value = -17000
$ws.Range("I9").Value = -0.07834101382488479
